{"js": "// Insert four new paragraphs into the document body, right after the\n// existing \"TBD: Sets. IDs. Hashing (bitstring), layout rules, set\n// encodings.\" paragraph (and before the trailing empty paragraph that\n// closes the body):\n//   1. an empty paragraph\n//   2. \"TBD: Ontology matching: Triadic FCA Context (object, attribute,\n//      condition). Fuzzy / rough sets (papers).\"\n//   3. an empty paragraph\n//   4. \"TBD: Protocol: Endpoints SPI / API implementation. Uniform\n//      Hypermedia interface. Messaging layer (levels). JDBC / JAF / JMS /\n//      JCA Connectors / Adapters.\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Find the anchor paragraph: the one whose text is the \"TBD: Sets...\" line.\nconst anchorText =\n  \"TBD: Sets. IDs. Hashing (bitstring), layout rules, set encodings.\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === anchorText) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  // Fallback: use the last non-empty paragraph before the trailing empty\n  // paragraph at the end of the body.\n  const items = paragraphs.items;\n  anchor = items[items.length - 2];\n}\n\nlet cursor = anchor.insertParagraph(\"\", Word.InsertLocation.after);\ncursor = cursor.insertParagraph(\n  \"TBD: Ontology matching: Triadic FCA Context (object, attribute, condition). Fuzzy / rough sets (papers).\",\n  Word.InsertLocation.after\n);\ncursor = cursor.insertParagraph(\"\", Word.InsertLocation.after);\ncursor = cursor.insertParagraph(\n  \"TBD: Protocol: Endpoints SPI / API implementation. Uniform Hypermedia interface. Messaging layer (levels). JDBC / JAF / JMS / JCA Connectors / Adapters.\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n", "ps1": "# Insert four new paragraphs into the document body, right after the\n# existing \"TBD: Sets. IDs. Hashing (bitstring), layout rules, set\n# encodings.\" paragraph (and before the trailing empty paragraph that\n# closes the body):\n#   1. an empty paragraph\n#   2. \"TBD: Ontology matching: Triadic FCA Context (object, attribute,\n#      condition). Fuzzy / rough sets (papers).\"\n#   3. an empty paragraph\n#   4. \"TBD: Protocol: Endpoints SPI / API implementation. Uniform\n#      Hypermedia interface. Messaging layer (levels). JDBC / JAF / JMS /\n#      JCA Connectors / Adapters.\"\n\n$d = $word.ActiveDocument\n\n$anchorText = \"TBD: Sets. IDs. Hashing (bitstring), layout rules, set encodings.\"\n\n$findRange = $d.Content\n$found = $findRange.Find.Execute($anchorText)\nif (-not $found) {\n    throw \"Could not find anchor paragraph '$anchorText'\"\n}\n\n# Work out the 1-based Paragraphs index of the anchor paragraph from the\n# number of paragraphs that precede the match.\n$anchorIndex = $d.Range(0, $findRange.Start).Paragraphs.Count + 1\n\n$newTexts = @(\n    \"\",\n    \"TBD: Ontology matching: Triadic FCA Context (object, attribute, condition). Fuzzy / rough sets (papers).\",\n    \"\",\n    \"TBD: Protocol: Endpoints SPI / API implementation. Uniform Hypermedia interface. Messaging layer (levels). JDBC / JAF / JMS / JCA Connectors / Adapters.\"\n)\n\n$currentIndex = $anchorIndex\nforeach ($text in $newTexts) {\n    $currentPara = $d.Paragraphs.Item($currentIndex)\n    $currentPara.Range.InsertParagraphAfter()\n    $currentIndex = $currentIndex + 1\n    if ($text -ne \"\") {\n        $d.Paragraphs.Item($currentIndex).Range.Text = $text\n    }\n}\n"}
